$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1 / index 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 47
$ws1.Range("F3").Value = 162
$ws1.Range("F5").Value = 13
$ws1.Range("F7").Value = 1609
$ws1.Range("F9").Value = 15
$ws1.Range("F10").Value = 1379
$ws1.Range("F11").Value = 120
$ws1.Range("F12").Value = 21
$ws1.Range("F13").Value = 238
$ws1.Range("F14").Value = 177
$ws1.Range("F17").Value = 10
$ws1.Range("F18").Value = 250
$ws1.Range("F19").Value = 135
$ws1.Range("F20").Value = 204
$ws1.Range("F21").Value = 191

# Sheet "全部类型" (sheet4 / index 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 47
$ws4.Range("F3").Value = 162
$ws4.Range("F5").Value = 13
$ws4.Range("F7").Value = 1609
$ws4.Range("F10").Value = 15
$ws4.Range("F11").Value = 1379
$ws4.Range("F12").Value = 120
$ws4.Range("F13").Value = 21
$ws4.Range("F14").Value = 238
$ws4.Range("F15").Value = 177
$ws4.Range("F18").Value = 10
$ws4.Range("F19").Value = 250
$ws4.Range("F20").Value = 135
$ws4.Range("F21").Value = 204
$ws4.Range("F22").Value = 191
